$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "永辉超市"
$ws.Range("A3").Value = "永辉超市"
$ws.Range("B3").Value = "永辉超市"
$ws.Range("C3").Value = "航天发展"
$ws.Range("A4").Value = "龙洲股份"
$ws.Range("B4").Value = "雷科防务"
$ws.Range("C4").Value = "雷科防务"
$ws.Range("A5").Value = "合富中国"
$ws.Range("B5").Value = "C摩尔-U"
$ws.Range("C5").Value = "实达集团"
$ws.Range("A6").Value = "航天机电"
$ws.Range("B6").Value = "达华智能"
$ws.Range("C6").Value = "龙洲股份"
$ws.Range("A7").Value = "达华智能"
$ws.Range("B7").Value = "龙洲股份"
$ws.Range("C7").Value = "海欣食品"
$ws.Range("A8").Value = "雷科防务"
$ws.Range("B8").Value = "实达集团"
$ws.Range("C8").Value = "合力泰"
$ws.Range("A9").Value = "航天动力"
$ws.Range("B9").Value = "合富中国"
$ws.Range("C9").Value = "安记食品"
$ws.Range("A10").Value = "实达集团"
$ws.Range("B10").Value = "航天机电"
$ws.Range("C10").Value = "航天动力"
$ws.Range("A11").Value = "安记食品"
$ws.Range("B11").Value = "工业富联"
$ws.Range("C11").Value = "工业富联"
$ws.Range("A12").Value = "C摩尔-U"
$ws.Range("B12").Value = "航天动力"
$ws.Range("C12").Value = "平潭发展"
$ws.Range("A13").Value = "乾照光电"
$ws.Range("B13").Value = "安妮股份"
$ws.Range("C13").Value = "达华智能"
$ws.Range("A14").Value = "工业富联"
$ws.Range("B14").Value = "航天科技"
$ws.Range("C14").Value = "合富中国"
$ws.Range("A15").Value = "中际旭创"
$ws.Range("B15").Value = "胜宏科技"
$ws.Range("C15").Value = "东百集团"
$ws.Range("A16").Value = "安妮股份"
$ws.Range("B16").Value = "中际旭创"
$ws.Range("A17").Value = "海欣食品"
$ws.Range("B17").Value = "合力泰"
$ws.Range("C17").Value = "安妮股份"
$ws.Range("A18").Value = "平潭发展"
$ws.Range("B18").Value = "平潭发展"
$ws.Range("C18").Value = "海王生物"
$ws.Range("A19").Value = "合力泰"
$ws.Range("B19").Value = "乾照光电"
$ws.Range("A20").Value = "鸿博股份"
$ws.Range("B20").Value = "东百集团"
$ws.Range("C20").Value = "厦门港务"
$ws.Range("A21").Value = "胜宏科技"
$ws.Range("B21").Value = "海欣食品"
$ws.Range("C21").Value = "博纳影业"
